$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Materials")

# Delete the suborder/infraorder/superfamily columns (AR:AT), then the
# Taxon_Local_ID column (A), from the Materials sheet. Deleting AR:AT first
# (while column A is still in place) keeps the AR:AT addresses correct;
# deleting A afterwards then shifts everything else left by one.
$ws.Range("AR:AT").EntireColumn.Delete()
$ws.Range("A:A").EntireColumn.Delete()

# The scientificNameAuthorship row's template value changes from
# ${summary.Author} to ${summary.authority}. Locate it with Find rather than
# a hardcoded address since the column deletions above already shifted it.
$found = $ws.Cells.Find("`${summary.Author}")
if ($found -ne $null) {
    $found.Value = "`${summary.authority}"
}
